$d = $word.ActiveDocument

# Locate template ranges to clone run formatting from (preserves exact rPr, including
# the presence/absence of an explicit w:lang element).
$tmplScope1 = $d.Range(0, $d.Content.End)
$tmplScope1.Find.Execute("RFID ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$noLangFT = $tmplScope1.FormattedText

$tmplScope2 = $d.Range(0, $d.Content.End)
$tmplScope2.Find.Execute("Фигура 7. Индуктивно-резонансен пренос на енергия ……………….. 16", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bgFT = $tmplScope2.FormattedText

function Add-FigRun {
    param($para, [string]$text, [bool]$noLang)
    $textEnd = $para.Range.End - 1
    $ip = $d.Range($textEnd, $textEnd)
    if ($noLang) {
        $ip.FormattedText = $noLangFT
    } else {
        $ip.FormattedText = $bgFT
    }
    $newEnd = $para.Range.End - 1
    $target = $d.Range($textEnd, $newEnd)
    $target.Text = $text
}

function New-FigParagraph {
    $lastPara = $d.Paragraphs.Last
    $r = $lastPara.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    return $d.Paragraphs.Last
}

# --- Фигура 8 ---
$p8 = New-FigParagraph
Add-FigRun $p8 "Фигура 8." $false
Add-FigRun $p8 " " $true
Add-FigRun $p8 "Снимка на" $false
Add-FigRun $p8 " " $false
Add-FigRun $p8 "Arduino Uno REV3" $true
Add-FigRun $p8 "………………………………" $false

# --- Фигура 9 ---
$p9 = New-FigParagraph
Add-FigRun $p9 "Фигура 9. Снимка на " $false
Add-FigRun $p9 "MRFC522" $true
Add-FigRun $p9 " " $false
Add-FigRun $p9 "RFID" $true
Add-FigRun $p9 " четец" $false
Add-FigRun $p9 "………………………….." $false

# --- Фигура 10 ---
$p10 = New-FigParagraph
Add-FigRun $p10 " " $false
Add-FigRun $p10 "Фигура 10" $false
Add-FigRun $p10 ". Снимка на " $false
Add-FigRun $p10 "MIFARE 1K Classic " $true
Add-FigRun $p10 "ключодържател ……" $false
Add-FigRun $p10 "…" $true
Add-FigRun $p10 ".." $true

# --- Фигура 11 ---
$p11 = New-FigParagraph
Add-FigRun $p11 " Фигура 11" $false
Add-FigRun $p11 ". Снимка на " $false
Add-FigRun $p11 "ESP-01 ESP8266 " $true
Add-FigRun $p11 "адаптер……………………" $false
Add-FigRun $p11 "." $false
Add-FigRun $p11 "." $false

# --- Фигура 12 ---
$p12 = New-FigParagraph
Add-FigRun $p12 " " $true
Add-FigRun $p12 "Фигура 12. Снимка на " $false
Add-FigRun $p12 "LED" $true
Add-FigRun $p12 "……………………………………………" $false
